$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert the two new columns (this shifts the data table in rows 5-7 and column widths correctly)
$ws.Columns("H:H").Insert()
$ws.Columns("AC:AC").Insert()

# 2) Set header values for the two newly inserted columns
$ws.Range("H5").Value = "Supplier ID"
$ws.Range("AC5").Value = "Approved By"

# 3) Fix up column widths for the two new columns to match the target widths
$ws.Columns("H:H").ColumnWidth = 13.166666666666666
$ws.Columns("AC:AC").ColumnWidth = 11.451822916666666

# 4) Fix up the "IRN REPORT" merged header box: it moves from G1:I2 to H1:J2
$ws.Range("G1:I2").UnMerge()
$ws.Range("G1:G2").Clear()
$ws.Range("H1:J2").Merge()
$ws.Range("H1").Value = "IRN REPORT"

# 5) Update the selection / view
$ws.Range("AC6").Select()
